{"js": "// 1. Change the title text \"LAB - 2\" -> \"LAB - 3\" in the first paragraph.\n// 2. Relocate the \"_GoBack\" bookmark from around the last paragraph\n//    (\"UPDATE STU SET Sem2_marks ...\") to a collapsed position right\n//    after the new \"LAB - 3\" text in the first paragraph.\n\n// --- Step 1: LAB - 2 -> LAB - 3 --------------------------------------\nconst results = context.document.body.search(\"LAB - 2\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"LAB - 3\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Step 2: drop the old _GoBack bookmark ---------------------------\nconst oldBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldBookmark.load(\"isNullObject\");\nawait context.sync();\n\nif (!oldBookmark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Step 3: add a fresh, collapsed _GoBack bookmark right after the\n// \"LAB - 3\" run in paragraph 1 (before its paragraph mark).\nconst firstParagraph = context.document.body.paragraphs.getFirst();\nconst endOfFirstParagraph = firstParagraph.getRange(\"End\");\nendOfFirstParagraph.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop script\n# 1. Change the title text \"LAB - 2\" -> \"LAB - 3\" in the first paragraph.\n# 2. Relocate the \"_GoBack\" bookmark from around the last paragraph\n#    (\"UPDATE STU SET Sem2_marks ...\") to a collapsed position right\n#    after the new \"LAB - 3\" text in the first paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: LAB - 2 -> LAB - 3 --------------------------------------\n$find = $d.Content.Find\n$find.Execute(\"LAB - 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"LAB - 3\", 2)\n\n# --- Step 2: drop the old _GoBack bookmark ---------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Step 3: add a fresh, collapsed _GoBack bookmark right after the\n# \"LAB - 3\" run in paragraph 1 (before its paragraph mark). Adding a\n# bookmark directly at a degenerate paragraph-boundary position can\n# make it swallow the whole paragraph, so instead we insert a\n# throwaway marker character, wrap the bookmark tightly around it\n# (a real, non-degenerate range), then delete the marker again -\n# leaving the bookmark collapsed exactly where we want it.\n$p1 = $d.Paragraphs.Item(1)\n$pos = $p1.Range.End - 1\n$marker = $d.Range($pos, $pos)\n$marker.InsertAfter(\"X\")\n$bmRange = $d.Range($pos, $pos + 1)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n$cleanup = $d.Range($pos, $pos + 1)\n$cleanup.Text = \"\"\n"}
